$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 becomes 7800
$ws.Range("C2").Value = 7800

# C3:C68 all become 7310
$ws.Range("C3:C68").Value = 7310
